$wb = $excel.ActiveWorkbook

# Overview sheet: G2 "Latest HO Xliff Generate Date" for d1ac55cc...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 21:04:03"

# zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 21:03:57"
$wsZhCn.Range("K2").Value = "2016-08-27 21:04:26"

# de-de sheet: K2 "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 21:04:33"
